$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5867746666666668
$ws.Range("H2").Value = 1.760324
$ws.Range("I2").Value = 0.5257388407083505
$ws.Range("J2").Value = 0.5257388407083505
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.749051
$ws.Range("N2").Value = 5.247153000000001
$ws.Range("O2").Value = 0.003644723415756578
$ws.Range("P2").Value = 0.003644723415756579
$ws.Range("Q2").Value = 1.026298817508
$ws.Range("R2").Value = 9.236689357572002
$ws.Range("S2").Value = 0.001916172663302443
$ws.Range("T2").Value = 0.001916172663302443
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5867746666666668
$ws.Range("H3").Value = 1.760324
$ws.Range("I3").Value = 0.5257388407083505
$ws.Range("J3").Value = 0.5257388407083505
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 143.0355533333334
$ws.Range("N3").Value = 429.10666
$ws.Range("O3").Value = 0.2980616520156925
$ws.Range("P3").Value = 0.2980616520156925
$ws.Range("Q3").Value = 83.92963912864892
$ws.Range("R3").Value = 755.3667521578402
$ws.Range("S3").Value = 0.156702587390346
$ws.Range("T3").Value = 0.156702587390346
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5867746666666668
$ws.Range("H4").Value = 1.760324
$ws.Range("I4").Value = 0.5257388407083505
$ws.Range("J4").Value = 0.5257388407083505
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 132.804812
$ws.Range("N4").Value = 398.414436
$ws.Range("O4").Value = 0.2767425352500014
$ws.Range("P4").Value = 0.2767425352500014
$ws.Range("Q4").Value = 77.92649929302935
$ws.Range("R4").Value = 701.3384936372642
$ws.Range("S4").Value = 0.1454942996570255
$ws.Range("T4").Value = 0.1454942996570256
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.5867746666666668
$ws.Range("H5").Value = 1.760324
$ws.Range("I5").Value = 0.5257388407083505
$ws.Range("J5").Value = 0.5257388407083505
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 125.707184
$ws.Range("N5").Value = 377.121552
$ws.Range("O5").Value = 0.261952291301752
$ws.Range("P5").Value = 0.261952291301752
$ws.Range("Q5").Value = 73.76179098920535
$ws.Range("R5").Value = 663.8561189028482
$ws.Range("S5").Value = 0.1377184939498792
$ws.Range("T5").Value = 0.1377184939498792
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.5867746666666668
$ws.Range("H6").Value = 1.760324
$ws.Range("I6").Value = 0.5257388407083505
$ws.Range("J6").Value = 0.5257388407083505
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 53.58648166666666
$ws.Range("N6").Value = 160.759445
$ws.Range("O6").Value = 0.111665071229204
$ws.Range("P6").Value = 0.1116650712292041
$ws.Range("Q6").Value = 31.44318991779778
$ws.Range("R6").Value = 282.98870926018
$ws.Range("S6").Value = 0.05870666509565713
$ws.Range("T6").Value = 0.05870666509565713
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.5867746666666668
$ws.Range("H7").Value = 1.760324
$ws.Range("I7").Value = 0.5257388407083505
$ws.Range("J7").Value = 0.5257388407083505
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.00271466666667
$ws.Range("N7").Value = 69.008144
$ws.Range("O7").Value = 0.04793372678759353
$ws.Range("P7").Value = 0.04793372678759355
$ws.Range("Q7").Value = 13.49741023096178
$ws.Range("R7").Value = 121.476692078656
$ws.Range("S7").Value = 0.02520062195214023
$ws.Range("T7").Value = 0.02520062195214024
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.5293206666666667
$ws.Range("H8").Value = 1.587962
$ws.Range("I8").Value = 0.4742611592916495
$ws.Range("J8").Value = 0.4742611592916495
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.749051
$ws.Range("N8").Value = 5.247153000000001
$ws.Range("O8").Value = 0.003644723415756578
$ws.Range("P8").Value = 0.003644723415756579
$ws.Range("Q8").Value = 0.9258088413540002
$ws.Range("R8").Value = 8.332279572186001
$ws.Range("S8").Value = 0.001728550752454135
$ws.Range("T8").Value = 0.001728550752454136
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.5293206666666667
$ws.Range("H9").Value = 1.587962
$ws.Range("I9").Value = 0.4742611592916495
$ws.Range("J9").Value = 0.4742611592916495
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 143.0355533333334
$ws.Range("N9").Value = 429.10666
$ws.Range("O9").Value = 0.2980616520156925
$ws.Range("P9").Value = 0.2980616520156925
$ws.Range("Q9").Value = 75.71167444743557
$ws.Range("R9").Value = 681.40507002692
$ws.Range("S9").Value = 0.1413590646253465
$ws.Range("T9").Value = 0.1413590646253465
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.5293206666666667
$ws.Range("H10").Value = 1.587962
$ws.Range("I10").Value = 0.4742611592916495
$ws.Range("J10").Value = 0.4742611592916495
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 132.804812
$ws.Range("N10").Value = 398.414436
$ws.Range("O10").Value = 0.2767425352500014
$ws.Range("P10").Value = 0.2767425352500014
$ws.Range("Q10").Value = 70.29633162438134
$ws.Range("R10").Value = 632.6669846194321
$ws.Range("S10").Value = 0.1312482355929758
$ws.Range("T10").Value = 0.1312482355929759
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.5293206666666667
$ws.Range("H11").Value = 1.587962
$ws.Range("I11").Value = 0.4742611592916495
$ws.Range("J11").Value = 0.4742611592916495
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 125.707184
$ws.Range("N11").Value = 377.121552
$ws.Range("O11").Value = 0.261952291301752
$ws.Range("P11").Value = 0.261952291301752
$ws.Range("Q11").Value = 66.53941043966933
$ws.Range("R11").Value = 598.8546939570241
$ws.Range("S11").Value = 0.1242337973518728
$ws.Range("T11").Value = 0.1242337973518728
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.5293206666666667
$ws.Range("H12").Value = 1.587962
$ws.Range("I12").Value = 0.4742611592916495
$ws.Range("J12").Value = 0.4742611592916495
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 53.58648166666666
$ws.Range("N12").Value = 160.759445
$ws.Range("O12").Value = 0.111665071229204
$ws.Range("P12").Value = 0.1116650712292041
$ws.Range("Q12").Value = 28.36443220012111
$ws.Range("R12").Value = 255.27988980109
$ws.Range("S12").Value = 0.05295840613354692
$ws.Range("T12").Value = 0.05295840613354694
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.5293206666666667
$ws.Range("H13").Value = 1.587962
$ws.Range("I13").Value = 0.4742611592916495
$ws.Range("J13").Value = 0.4742611592916495
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 23.00271466666667
$ws.Range("N13").Value = 69.008144
$ws.Range("O13").Value = 0.04793372678759353
$ws.Range("P13").Value = 0.04793372678759355
$ws.Range("Q13").Value = 12.17581226250311
$ws.Range("R13").Value = 109.582310362528
$ws.Range("S13").Value = 0.0227331048354533
$ws.Range("T13").Value = 0.02273310483545331

Write-Output "done"